$d = $word.ActiveDocument

# 1. "Imajući to u obziru u podvodnim" -> "Imajući to u vidu u podvodnim"  (obziru -> vidu)
$d.Content.Find.Execute("Imajući to u obziru u podvodnim", $true, $false, $false, $false, $false, $true, 1, $false, "Imajući to u vidu u podvodnim", 2)

# 2. "koje si radile" -> "koje su radile"
$d.Content.Find.Execute("koje si radile", $true, $false, $false, $false, $false, $true, 1, $false, "koje su radile", 2)

# 3. "Takođe u ISP se radi projekat" -> "Takođe u ISP se radio projekat"
$d.Content.Find.Execute("Takođe u ISP se radi projekat", $true, $false, $false, $false, $false, $true, 1, $false, "Takođe u ISP se radio projekat", 2)

# 4. "carrief frekvenciju" -> "carrier frekvenciju"
$d.Content.Find.Execute("carrief frekvenciju", $true, $false, $false, $false, $false, $true, 1, $false, "carrier frekvenciju", 2)

# 5. "1.1 Postava ekssperimenta" -> "1.1 Postava eksperimenta"
$d.Content.Find.Execute("1.1 Postava ekssperimenta", $true, $false, $false, $false, $false, $true, 1, $false, "1.1 Postava eksperimenta", 2)

# 6. "kako ne bi dolazilo da aliasinga" -> "kako ne bi dolazilo do aliasinga"
$d.Content.Find.Execute("kako ne bi dolazilo da aliasinga", $true, $false, $false, $false, $false, $true, 1, $false, "kako ne bi dolazilo do aliasinga", 2)

# 6b. Make the word "aliasinga" italic
$r = $d.Content
$r.Find.Execute("aliasinga", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Italic = 1

# 7. "Ovi vaktori mogu" -> "Ovi faktori mogu"
$d.Content.Find.Execute("Ovi vaktori mogu", $true, $false, $false, $false, $false, $true, 1, $false, "Ovi faktori mogu", 2)

# 8. "eksperimenta Br. 1. jr to" -> "eksperimenta Br. 1. je to"
$d.Content.Find.Execute("eksperimenta Br. 1. jr to", $true, $false, $false, $false, $false, $true, 1, $false, "eksperimenta Br. 1. je to", 2)
